$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.215.87"
$ws.Range("E2").Value = "  +1.33%  "
$ws.Range("D3").Value = "'3.080.70"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'560.26"
$ws.Range("E5").Value = "  +1.95%  "
$ws.Range("D6").Value = "'144.98"
$ws.Range("E6").Value = "  +3.57%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'3.078.06"
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("E10").Value = "  +1.81%  "
$ws.Range("D11").Value = "'6.17"
$ws.Range("E11").Value = "  -3.58%  "
$ws.Range("D12").Value = "'0.471"
$ws.Range("E12").Value = "  +4.05%  "
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("D14").Value = "'35.15"
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("D15").Value = "'3.577.18"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").Value = "'64.270.00"
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("D17").Value = "'3.076.09"
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("D19").Value = "'6.76"
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("D20").Value = "'478.83"
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("D21").Value = "'13.91"
$ws.Range("E21").Value = "  +1.80%  "
$ws.Range("D22").Value = "'0.675"
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").Value = "'7.56"
$ws.Range("E23").Value = "  +5.12%  "
$ws.Range("D24").Value = "'13.76"
$ws.Range("E24").Value = "  +10.20%  "
$ws.Range("D25").Value = "'81.15"
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'2.81"
$ws.Range("E27").Value = "  +2.25%  "
$ws.Range("D28").Value = "'8.04"
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("D29").Value = "'2.08"
$ws.Range("E29").Value = "  +4.56%  "
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("D31").Value = "'26.17"
$ws.Range("E31").Value = "  +0.71%  "
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("D33").Value = "'2.50"
$ws.Range("E33").Value = "  +3.07%  "
$ws.Range("D34").Value = "'5.57"
$ws.Range("E34").Value = "  -1.98%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'6.20"
$ws.Range("E35").Value = "  +3.91%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "'55.59"
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").Value = "'457.29"
$ws.Range("E37").Value = "  -1.35%  "
$ws.Range("E38").Value = "  +16.91%  "
$ws.Range("D39").Value = "'0.0406"
$ws.Range("E39").Value = "  +2.69%  "
$ws.Range("D40").Value = "'0.0825"
$ws.Range("E40").Value = "  +1.17%  "
$ws.Range("D41").Value = "'2.962.99"
$ws.Range("E41").Value = "  -3.50%  "
$ws.Range("D42").Value = "'8.25"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("E43").Value = "  -1.98%  "
$ws.Range("D44").Value = "'27.91"
$ws.Range("E44").Value = "  -1.36%  "
$ws.Range("D45").Value = "'0.261"
$ws.Range("E45").Value = "  +3.34%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.14"
$ws.Range("E46").Value = "  +4.60%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "'1.00"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "'0.112"
$ws.Range("E48").Value = "  +2.28%  "
$ws.Range("D49").Value = "'121.10"
$ws.Range("E49").Value = "  +3.72%  "
$ws.Range("E50").Value = "  +1.10%  "
$ws.Range("D51").Value = "'2.08"
$ws.Range("E51").Value = "  +0.54%  "
